$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is numeric-looking must keep Text format
# so Excel does not silently convert them to Number (losing trailing zeros / exact text).

$ws.Range('D2').Value = '26.727.22'
$ws.Range('E2').Value = '  +0.48%  '

$ws.Range('D3').Value = '1.648.57'
$ws.Range('E3').Value = '  +1.13%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.03'
$ws.Range('E5').Value = '  +1.09%  '

$ws.Range('E6').Value = '  +2.40%  '

$ws.Range('E7').Value = '  +0.16%  '

$ws.Range('E8').Value = '  -0.17%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0628'
$ws.Range('E9').Value = '  +0.83%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.24'
$ws.Range('E10').Value = '  +1.80%  '

$ws.Range('D12').Value = '1.879.34'
$ws.Range('E12').Value = '  +1.13%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.15'
$ws.Range('E13').Value = '  +2.05%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.606.71'
$ws.Range('E14').Value = '  -1.07%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.531'
$ws.Range('E15').Value = '  +1.50%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.19'
$ws.Range('E16').Value = '  +5.22%  '

$ws.Range('D17').Value = '26.793.67'
$ws.Range('E17').Value = '  +0.77%  '

$ws.Range('D18').Value = '0.0₃0747'
$ws.Range('E18').Value = '  +1.28%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '218.31'
$ws.Range('E19').Value = '  +4.61%  '

$ws.Range('E20').Value = '  +0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.36'
$ws.Range('E21').Value = '  +1.92%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.38'
$ws.Range('E22').Value = '  +3.44%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.46'
$ws.Range('E23').Value = '  +0.45%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.60'
$ws.Range('E25').Value = '  +0.66%  '

$ws.Range('E26').Value = '  +0.15%  '

$ws.Range('E27').Value = '  +0.29%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.95'
$ws.Range('E28').Value = '  +1.31%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.80'
$ws.Range('E29').Value = '  +2.93%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0522'
$ws.Range('E30').Value = '  +0.65%  '

$ws.Range('E31').Value = '  +0.49%  '

$ws.Range('E32').Value = '  +4.29%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.00'
$ws.Range('E33').Value = '  +2.58%  '

$ws.Range('D34').Value = '1.281.64'
$ws.Range('E34').Value = '  +10.10%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').Value = '  +2.03%  '

$ws.Range('E36').Value = '  +1.39%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0178'
$ws.Range('E37').Value = '  +4.42%  '

$ws.Range('E38').Value = '  +2.80%  '

$ws.Range('E39').Value = '  +0.42%  '

$ws.Range('E40').Value = '  +0.07%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.28'
$ws.Range('E41').Value = '  -1.37%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.808'
$ws.Range('E42').Value = '  +2.37%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.39'
$ws.Range('E43').Value = '  +0.59%  '

$ws.Range('D44').Value = '1.788.95'
$ws.Range('E44').Value = '  +1.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.52'
$ws.Range('E45').Value = '  +1.17%  '

$ws.Range('E46').Value = '  +4.63%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.07'
$ws.Range('E47').Value = '  +2.90%  '

$ws.Range('E48').Value = '  -2.12%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0514'
$ws.Range('E49').Value = '  +0.24%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.72'
$ws.Range('E50').Value = '  +2.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0967'
$ws.Range('E51').Value = '  +3.14%  '
